# The synthetic array gains a new "statut_name" column right after
# "statut_label" (i.e. before the former column C "NCTId"). All the
# columns from the old C ("NCTId") onward shift one position to the
# right (new D..M), and the new column C is populated with a human
# readable label describing the "statut" status of each trial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C ("NCTId"), pushing
# NCTId..intervention_type from C..L to D..M.
$ws.Range("C1").EntireColumn.Insert()

# Populate the newly inserted column with its header and values.
$ws.Range("C1").Value = "statut_name"
$ws.Range("C2").Value = "résultat et / ou publication posté dans les 12 mois"
$ws.Range("C3").Value = "pas de résultat ni de publication"
$ws.Range("C4").Value = "pas de résultat ni de publication"
